$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Template")

# Populate the new row 7 values (Subgroup, Name, Filename, Description/Notes)
$ws.Range("B7").Value = "Animations"
$ws.Range("C7").Value = "MageAnimations"
$ws.Range("D7").Value = "Frank_Mage"
$ws.Range("E7").Value = "Animation clips in Mage animation folder are copied from here due to import errors"

# Update the active selection to match the new cursor position
$ws.Range("D10").Select()
